$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's portfolio row (row 31) below the existing data (rows 1-30).
$row = 31

# Force the date column to be stored as plain text (matching the existing
# "Date" column cells, which are text strings like "2025-08-17"), then
# clear the temporary formatting so the new cell doesn't carry an explicit
# number-format style, consistent with the other data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-15"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 58.06999969482422
$ws.Cells.Item($row, 3).Value = 712.9000244140625
$ws.Cells.Item($row, 4).Value = 323.2999877929688
